$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.791024514463149
$ws.Range("D2").Value = 10.01111283468003
$ws.Range("E2").Value = 14.00626619514807
$ws.Range("F2").Value = 28.77672641412201
$ws.Range("G2").Value = 27.65552021426359
$ws.Range("H2").Value = 14.00198121041143
$ws.Range("J2").Value = 9.762207578314579
$ws.Range("K2").Value = 10.35965400867802
$ws.Range("N2").Value = 18.02549802065435
$ws.Range("O2").Value = 21.16853553241575
$ws.Range("B3").Value = 7.718250230697997
$ws.Range("D3").Value = 9.955422248793003
$ws.Range("E3").Value = 13.94238617244795
$ws.Range("F3").Value = 28.7915958502273
$ws.Range("G3").Value = 27.66304705817157
$ws.Range("H3").Value = 14.04237902774122
$ws.Range("J3").Value = 9.766875411747282
$ws.Range("K3").Value = 9.976793225029994
$ws.Range("N3").Value = 18.0720333991718
$ws.Range("O3").Value = 21.22380505608925
$ws.Range("B4").Value = 7.675042152190723
$ws.Range("D4").Value = 9.922859528833959
$ws.Range("E4").Value = 13.90595327354747
$ws.Range("F4").Value = 28.80788573680208
$ws.Range("G4").Value = 27.67677249694476
$ws.Range("H4").Value = 14.06944064302012
$ws.Range("J4").Value = 9.771286968313774
$ws.Range("K4").Value = 9.734595137235482
$ws.Range("N4").Value = 18.10238535455759
$ws.Range("O4").Value = 21.26237707074541
$ws.Range("B5").Value = 7.657825488934114
$ws.Range("D5").Value = 9.910011044993201
$ws.Range("E5").Value = 13.8918200475837
$ws.Range("F5").Value = 28.81632367889824
$ws.Range("G5").Value = 27.68465135212408
$ws.Range("H5").Value = 14.08103585604695
$ws.Range("J5").Value = 9.773473835331908
$ws.Range("K5").Value = 9.634246076997725
$ws.Range("N5").Value = 18.11520225271807
$ws.Range("O5").Value = 21.27925897922249
$ws.Range("B6").Value = 7.654990857280001
$ws.Range("D6").Value = 9.907903301509888
$ws.Range("E6").Value = 13.88951666007131
$ws.Range("F6").Value = 28.81783344828334
$ws.Range("G6").Value = 27.68609753534222
$ws.Range("H6").Value = 14.08299549478396
$ws.Range("J6").Value = 9.773860477172255
$ws.Range("K6").Value = 9.617487806571647
$ws.Range("N6").Value = 18.11735758680478
$ws.Range("O6").Value = 21.28213240794456
$ws.Range("B7").Value = 7.674808353396602
$ws.Range("D7").Value = 9.922684530847411
$ws.Range("E7").Value = 13.90575976394356
$ws.Range("F7").Value = 28.80799224871883
$ws.Range("G7").Value = 27.67686950565821
$ws.Range("H7").Value = 14.0695947230823
$ws.Range("J7").Value = 9.771314885056015
$ws.Range("K7").Value = 9.733248286509019
$ws.Range("N7").Value = 18.10255639169316
$ws.Range("O7").Value = 21.26260003833392
$ws.Range("B8").Value = 7.765637840171776
$ws.Range("D8").Value = 9.991578770384017
$ws.Range("E8").Value = 13.98366883388876
$ws.Range("F8").Value = 28.78036716093229
$ws.Range("G8").Value = 27.65622430177129
$ws.Range("H8").Value = 14.01544177963327
$ws.Range("J8").Value = 9.763496611136603
$ws.Range("K8").Value = 10.22920023321918
$ws.Range("N8").Value = 18.04117470299423
$ws.Range("O8").Value = 21.18662884634522
$ws.Range("B9").Value = 7.954555439031989
$ws.Range("D9").Value = 10.13912369891071
$ws.Range("E9").Value = 14.15801448204709
$ws.Range("F9").Value = 28.78300097761393
$ws.Range("G9").Value = 27.68805621424058
$ws.Range("H9").Value = 13.92716664708085
$ws.Range("J9").Value = 9.760402219964904
$ws.Range("K9").Value = 11.13989215454348
$ws.Range("N9").Value = 17.93488200431752
$ws.Range("O9").Value = 21.07453452363423
$ws.Range("B10").Value = 8.098687366860807
$ws.Range("D10").Value = 10.25440416961667
$ws.Range("E10").Value = 14.2984342341412
$ws.Range("F10").Value = 28.81950342726159
$ws.Range("G10").Value = 27.75552247163033
$ws.Range("H10").Value = 13.87324612950033
$ws.Range("J10").Value = 9.76554906507276
$ws.Range("K10").Value = 11.76502725463668
$ws.Range("N10").Value = 17.86531605557566
$ws.Range("O10").Value = 21.01478889754749
$ws.Range("B11").Value = 8.165142715843395
$ws.Range("D11").Value = 10.30818223547524
$ws.Range("E11").Value = 14.36480085170023
$ws.Range("F11").Value = 28.84358060397141
$ws.Range("G11").Value = 27.79574422357594
$ws.Range("H11").Value = 13.85109298595225
$ws.Range("J11").Value = 9.769490760150321
$ws.Range("K11").Value = 12.03880898853374
$ws.Range("N11").Value = 17.83550905639827
$ws.Range("O11").Value = 20.99254099665347
$ws.Range("B12").Value = 8.190412520538402
$ws.Range("D12").Value = 10.32872475162348
$ws.Range("E12").Value = 14.39027293410531
$ws.Range("F12").Value = 28.85376770863856
$ws.Range("G12").Value = 27.81233892006992
$ws.Range("H12").Value = 13.84304594781106
$ws.Range("J12").Value = 9.77121233643193
$ws.Range("K12").Value = 12.14088558157565
$ws.Range("N12").Value = 17.82448552724511
$ws.Range("O12").Value = 20.98482676799748
$ws.Range("B13").Value = 8.184965932339175
$ws.Range("D13").Value = 10.32429286943926
$ws.Range("E13").Value = 14.3847722058916
$ws.Range("F13").Value = 28.8515262516046
$ws.Range("G13").Value = 27.80870444108133
$ws.Range("H13").Value = 13.84476381230063
$ws.Range("J13").Value = 9.770831402275219
$ws.Range("K13").Value = 12.11897374795828
$ws.Range("N13").Value = 17.82684792274143
$ws.Range("O13").Value = 20.98645654050669
$ws.Range("B14").Value = 8.167219714121297
$ws.Range("D14").Value = 10.30986879722632
$ws.Range("E14").Value = 14.36688973331279
$ws.Range("F14").Value = 28.84439729853775
$ws.Range("G14").Value = 27.7970821977518
$ws.Range("H14").Value = 13.85042409808466
$ws.Range("J14").Value = 9.769627814214346
$ws.Range("K14").Value = 12.04723932392633
$ws.Range("N14").Value = 17.83459686232023
$ws.Range("O14").Value = 20.99189209361041
$ws.Range("B15").Value = 8.15636257804878
$ws.Range("D15").Value = 10.30105638075053
$ws.Range("E15").Value = 14.3559800164898
$ws.Range("F15").Value = 28.84016973958553
$ws.Range("G15").Value = 27.79014057760082
$ws.Range("H15").Value = 13.85393571669078
$ws.Range("J15").Value = 9.768920359205218
$ws.Range("K15").Value = 12.00308967685316
$ws.Range("N15").Value = 17.83937763718019
$ws.Range("O15").Value = 20.99531410540339
$ws.Range("B16").Value = 8.094360254907357
$ws.Range("D16").Value = 10.25091532737448
$ws.Range("E16").Value = 14.29414562804854
$ws.Range("F16").Value = 28.81807989196418
$ws.Range("G16").Value = 27.75308512482438
$ws.Range("H16").Value = 13.87474171586252
$ws.Range("J16").Value = 9.765323570926016
$ws.Range("K16").Value = 11.74691512530133
$ws.Range("N16").Value = 17.86730095647094
$ws.Range("O16").Value = 21.01634220559104
$ws.Range("B17").Value = 8.05653498481794
$ws.Range("D17").Value = 10.22048765938105
$ws.Range("E17").Value = 14.25683751241794
$ws.Range("F17").Value = 28.80643913091988
$ws.Range("G17").Value = 27.7327896105754
$ws.Range("H17").Value = 13.88811416369499
$ws.Range("J17").Value = 9.763526146837352
$ws.Range("K17").Value = 11.58698998858465
$ws.Range("N17").Value = 17.88490149351039
$ws.Range("O17").Value = 21.03050630523277
$ws.Range("B18").Value = 8.034863957547664
$ws.Range("D18").Value = 10.20311313067484
$ws.Range("E18").Value = 14.23561422803415
$ws.Range("F18").Value = 28.80044733568558
$ws.Range("G18").Value = 27.72201381509351
$ws.Range("H18").Value = 13.89602922694574
$ws.Range("J18").Value = 9.76264302996862
$ws.Range("K18").Value = 11.49401179752829
$ws.Range("N18").Value = 17.8951979834529
$ws.Range("O18").Value = 21.03911714303776
$ws.Range("B19").Value = 8.027541839584112
$ws.Range("D19").Value = 10.19725259716684
$ws.Range("E19").Value = 14.22846933246399
$ws.Range("F19").Value = 28.79853960217764
$ws.Range("G19").Value = 27.71851966700142
$ws.Range("H19").Value = 13.89874752364828
$ws.Range("J19").Value = 9.762369941092908
$ws.Range("K19").Value = 11.46236280060199
$ws.Range("N19").Value = 17.89871395419948
$ws.Range("O19").Value = 21.04211227006695
$ws.Range("B20").Value = 8.060552913775142
$ws.Range("D20").Value = 10.22371372990261
$ws.Range("E20").Value = 14.26078479080017
$ws.Range("F20").Value = 28.80760551519002
$ws.Range("G20").Value = 27.73485724360412
$ws.Range("H20").Value = 13.88666750154225
$ws.Range("J20").Value = 9.763701895568879
$ws.Range("K20").Value = 11.60411768941346
$ws.Range("N20").Value = 17.88300997344759
$ws.Range("O20").Value = 21.0289504761158
$ws.Range("B21").Value = 8.17242955423796
$ws.Range("D21").Value = 10.31410078279365
$ws.Range("E21").Value = 14.37213315009372
$ws.Range("F21").Value = 28.84646226035989
$ws.Range("G21").Value = 27.80045899015904
$ws.Range("H21").Value = 13.84875225444963
$ws.Range("J21").Value = 9.769975133922278
$ws.Range("K21").Value = 12.06835339718652
$ws.Range("N21").Value = 17.83231365872114
$ws.Range("O21").Value = 20.99027624420922
$ws.Range("B22").Value = 8.24614645521963
$ws.Range("D22").Value = 10.37420427009307
$ws.Range("E22").Value = 14.44688186819485
$ws.Range("F22").Value = 28.87808882042834
$ws.Range("G22").Value = 27.85127695343517
$ws.Range("H22").Value = 13.82596527732157
$ws.Range("J22").Value = 9.775408861015489
$ws.Range("K22").Value = 12.36241240608656
$ws.Range("N22").Value = 17.80071762154445
$ws.Range("O22").Value = 20.969142674963
$ws.Range("B23").Value = 8.206755187227623
$ws.Range("D23").Value = 10.34203640666378
$ws.Range("E23").Value = 14.40681215815018
$ws.Range("F23").Value = 28.86064084058637
$ws.Range("G23").Value = 27.82343043726464
$ws.Range("H23").Value = 13.83794468766789
$ws.Range("J23").Value = 9.772387162727602
$ws.Range("K23").Value = 12.20634497144917
$ws.Range("N23").Value = 17.81744061750884
$ws.Range("O23").Value = 20.98004259462058
$ws.Range("B24").Value = 8.058736172016159
$ws.Range("D24").Value = 10.22225485247631
$ws.Range("E24").Value = 14.25899952206895
$ws.Range("F24").Value = 28.80707601008008
$ws.Range("G24").Value = 27.73391968654749
$ws.Range("H24").Value = 13.88732083022192
$ws.Range("J24").Value = 9.763621971480035
$ws.Range("K24").Value = 11.59637747416522
$ws.Range("N24").Value = 17.88386457606965
$ws.Range("O24").Value = 21.02965240968975
$ws.Range("B25").Value = 7.90241954422569
$ws.Range("D25").Value = 10.09794778903954
$ws.Range("E25").Value = 14.10862188535001
$ws.Range("F25").Value = 28.77620938101921
$ws.Range("G25").Value = 27.67169588403157
$ws.Range("H25").Value = 13.94912783871837
$ws.Range("J25").Value = 9.759932866096459
$ws.Range("K25").Value = 10.90084340329563
$ws.Range("N25").Value = 17.96213566221779
$ws.Range("O25").Value = 21.10089658301655
